# Applies updated FFXIV market-board figures (columns H-N: currentAveragePrice*,
# LevePrice*, LeveProfit*) across the leve-profit worksheets, per the scheduled
# market-data refresh.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40000
$ws.Range("J3").Value = 40000
$ws.Range("L3").Value = 40000
$ws.Range("N3").Value = -40228
$ws.Range("H20").Value = 14000
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 27000
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 27000
$ws.Range("M20").Value = -770
$ws.Range("N20").Value = -27460
$ws.Range("H34").Value = 4705
$ws.Range("I34").Value = 1044.7778
$ws.Range("J34").Value = 26666.334
$ws.Range("K34").Value = 1044.7778
$ws.Range("L34").Value = 26666.334
$ws.Range("M34").Value = -841.7778000000001
$ws.Range("N34").Value = -27072.334
$ws.Range("H35").Value = 14000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 27000
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 27000
$ws.Range("M35").Value = -621
$ws.Range("N35").Value = -27758
$ws.Range("H36").Value = 4705
$ws.Range("I36").Value = 1044.7778
$ws.Range("J36").Value = 26666.334
$ws.Range("K36").Value = 1044.7778
$ws.Range("L36").Value = 26666.334
$ws.Range("M36").Value = -329.7778000000001
$ws.Range("N36").Value = -28096.334
$ws.Range("H44").Value = 27000
$ws.Range("J44").Value = 27000
$ws.Range("L44").Value = 27000
$ws.Range("N44").Value = -27924
$ws.Range("H47").Value = 19960
$ws.Range("I47").Value = 12900
$ws.Range("J47").Value = 24666.666
$ws.Range("K47").Value = 12900
$ws.Range("L47").Value = 24666.666
$ws.Range("M47").Value = -11928
$ws.Range("N47").Value = -26610.666
$ws.Range("H93").Value = 44450
$ws.Range("J93").Value = 44450
$ws.Range("L93").Value = 44450
$ws.Range("N93").Value = -49442
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H105").Value = 33000
$ws.Range("J105").Value = 33000
$ws.Range("L105").Value = 33000
$ws.Range("N105").Value = -39988
$ws.Range("H125").Value = 893.06665
$ws.Range("I125").Value = 728
$ws.Range("J125").Value = 1081.7142
$ws.Range("K125").Value = 6552
$ws.Range("L125").Value = 9735.427799999999
$ws.Range("M125").Value = -4092
$ws.Range("N125").Value = -14655.4278

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 14095.454
$ws.Range("I28").Value = 4933.2856
$ws.Range("J28").Value = 30129.25
$ws.Range("K28").Value = 4933.2856
$ws.Range("L28").Value = 30129.25
$ws.Range("M28").Value = -4741.2856
$ws.Range("N28").Value = -30513.25
$ws.Range("H31").Value = 12733.875
$ws.Range("I31").Value = 5974.2
$ws.Range("J31").Value = 24000
$ws.Range("K31").Value = 5974.2
$ws.Range("L31").Value = 24000
$ws.Range("M31").Value = -5680.2
$ws.Range("N31").Value = -24588
$ws.Range("H32").Value = 951328.9399999999
$ws.Range("I32").Value = 1007041.5
$ws.Range("J32").Value = 26500
$ws.Range("K32").Value = 1007041.5
$ws.Range("L32").Value = 26500
$ws.Range("M32").Value = -1006754.5
$ws.Range("N32").Value = -27074
$ws.Range("H93").Value = 27500
$ws.Range("H99").Value = 14095.454
$ws.Range("I99").Value = 4933.2856
$ws.Range("J99").Value = 30129.25
$ws.Range("K99").Value = 4933.2856
$ws.Range("L99").Value = 30129.25
$ws.Range("M99").Value = -1938.2856
$ws.Range("N99").Value = -36119.25

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 29139.334
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 29139.334
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 29139.334
$ws.Range("N28").Value = -29727.334
$ws.Range("M28").ClearContents()
$ws.Range("H40").Value = 28000
$ws.Range("J40").Value = 28000
$ws.Range("L40").Value = 28000
$ws.Range("N40").Value = -28530
$ws.Range("H44").Value = 22050
$ws.Range("J44").Value = 22050
$ws.Range("L44").Value = 22050
$ws.Range("N44").Value = -23044
$ws.Range("H101").Value = 23000
$ws.Range("J101").Value = 23000
$ws.Range("L101").Value = 23000
$ws.Range("N101").Value = -29490

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 9560.888999999999
$ws.Range("J96").Value = 9560.888999999999
$ws.Range("L96").Value = 9560.888999999999
$ws.Range("N96").Value = -15052.889

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 3533.3333
$ws.Range("I14").Value = 3533.3333
$ws.Range("K14").Value = 10599.9999
$ws.Range("M14").Value = -10426.9999

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 8663.909
$ws.Range("I12").Value = 1471.4286
$ws.Range("J12").Value = 21250.75
$ws.Range("K12").Value = 1471.4286
$ws.Range("L12").Value = 21250.75
$ws.Range("M12").Value = -1331.4286
$ws.Range("N12").Value = -21530.75
$ws.Range("H94").Value = 27500
$ws.Range("J94").Value = 27500
$ws.Range("L94").Value = 27500
$ws.Range("N94").Value = -28852
$ws.Range("H98").Value = 25124.875
$ws.Range("J98").Value = 25124.875
$ws.Range("L98").Value = 25124.875
$ws.Range("N98").Value = -31114.875
$ws.Range("H99").Value = 4897.2
$ws.Range("I99").Value = 3774.6667
$ws.Range("J99").Value = 15000
$ws.Range("K99").Value = 3774.6667
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -1528.6667
$ws.Range("N99").Value = -19492

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 501.5
$ws.Range("I9").Value = 335.5
$ws.Range("J9").Value = 999.5
$ws.Range("K9").Value = 335.5
$ws.Range("L9").Value = 999.5
$ws.Range("M9").Value = -111.5
$ws.Range("N9").Value = -1447.5
$ws.Range("H57").Value = 9991.571
$ws.Range("I57").Value = 1680.3334
$ws.Range("J57").Value = 16225
$ws.Range("K57").Value = 1680.3334
$ws.Range("L57").Value = 16225
$ws.Range("M57").Value = -1114.3334
$ws.Range("N57").Value = -17357

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 7775.857
$ws.Range("I61").Value = 3637.75
$ws.Range("J61").Value = 13293.333
$ws.Range("K61").Value = 3637.75
$ws.Range("L61").Value = 13293.333
$ws.Range("M61").Value = -3345.75
$ws.Range("N61").Value = -13877.333
